# Update the newly added iAuthor TC credential/sample row values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ctrTv372"
$ws.Range("B2").Value = 231011251
$ws.Range("C2").Value = "xcdylvi50"
$ws.Range("D2").Value = "EFrk&48#"
$ws.Range("F2").Value = "pgUnTswa"
$ws.Range("G2").Value = "uzgK"
